$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Schedule")
$ws2 = $wb.Worksheets.Item("Detailed")

# --- Schedule sheet updates ---
# Schedule sheet row 2 update + new rows 3-5
$r = 2
$ws1.Cells.Item($r, 1).Value = 46071.08333333334
$ws1.Cells.Item($r, 2).Value = 46071.25
$ws1.Cells.Item($r, 3).Value = 4
$ws1.Cells.Item($r, 4).Value = 15.12
$ws1.Cells.Item($r, 5).Value = 562.965507
$ws1.Cells.Item($r, 6).Value = 37.23316845238095
$r = 3
$ws1.Cells.Item($r, 1).Value = 46071.3125
$ws1.Cells.Item($r, 1).NumberFormat = $ws1.Cells.Item(2, 1).NumberFormat
$ws1.Cells.Item($r, 2).Value = 46071.60416666666
$ws1.Cells.Item($r, 2).NumberFormat = $ws1.Cells.Item(2, 2).NumberFormat
$ws1.Cells.Item($r, 3).Value = 7
$ws1.Cells.Item($r, 4).Value = 26.46
$ws1.Cells.Item($r, 5).Value = 194.5356855
$ws1.Cells.Item($r, 6).Value = 7.35206672335601
$r = 4
$ws1.Cells.Item($r, 1).Value = 46072.0625
$ws1.Cells.Item($r, 1).NumberFormat = $ws1.Cells.Item(2, 1).NumberFormat
$ws1.Cells.Item($r, 2).Value = 46072.22916666666
$ws1.Cells.Item($r, 2).NumberFormat = $ws1.Cells.Item(2, 2).NumberFormat
$ws1.Cells.Item($r, 3).Value = 4
$ws1.Cells.Item($r, 4).Value = 15.12
$ws1.Cells.Item($r, 5).Value = 639.4456185
$ws1.Cells.Item($r, 6).Value = 42.29137688492064
$r = 5
$ws1.Cells.Item($r, 1).Value = 46072.33333333334
$ws1.Cells.Item($r, 1).NumberFormat = $ws1.Cells.Item(2, 1).NumberFormat
$ws1.Cells.Item($r, 2).Value = 46072.625
$ws1.Cells.Item($r, 2).NumberFormat = $ws1.Cells.Item(2, 2).NumberFormat
$ws1.Cells.Item($r, 3).Value = 7
$ws1.Cells.Item($r, 4).Value = 26.46
$ws1.Cells.Item($r, 5).Value = 397.6421767500001
$ws1.Cells.Item($r, 6).Value = 15.02804900793651

# --- Detailed sheet: update existing rows ---
$ws2.Cells.Item(6, 5).Value = "ON"
$ws2.Cells.Item(7, 5).Value = "ON"
$ws2.Cells.Item(8, 5).Value = "ON"
$ws2.Cells.Item(9, 5).Value = "ON"
$ws2.Cells.Item(14, 2).Value = 84.79000000000001
$ws2.Cells.Item(14, 5).Value = "OFF"
$ws2.Cells.Item(15, 2).Value = 107.78769
$ws2.Cells.Item(15, 5).Value = "OFF"
$ws2.Cells.Item(16, 2).Value = 94.74899000000001
$ws2.Cells.Item(16, 3).Value = "historical"
$ws2.Cells.Item(16, 5).Value = "OFF"
$ws2.Cells.Item(17, 2).Value = 35.88
$ws2.Cells.Item(17, 3).Value = "historical"
$ws2.Cells.Item(18, 2).Value = 13.26863
$ws2.Cells.Item(18, 3).Value = "historical"
$ws2.Cells.Item(19, 3).Value = "historical"
$ws2.Cells.Item(20, 2).Value = 0.51
$ws2.Cells.Item(20, 3).Value = "historical"
$ws2.Cells.Item(21, 3).Value = "historical"
$ws2.Cells.Item(22, 2).Value = 0.51
$ws2.Cells.Item(22, 3).Value = "historical"
$ws2.Cells.Item(23, 2).Value = 0.51
$ws2.Cells.Item(23, 3).Value = "historical"
$ws2.Cells.Item(24, 3).Value = "historical"
$ws2.Cells.Item(25, 3).Value = "historical"
$ws2.Cells.Item(26, 2).Value = 22.07
$ws2.Cells.Item(26, 3).Value = "historical"
$ws2.Cells.Item(27, 3).Value = "historical"
$ws2.Cells.Item(28, 2).Value = 30.34515
$ws2.Cells.Item(28, 3).Value = "historical"
$ws2.Cells.Item(29, 3).Value = "historical"
$ws2.Cells.Item(30, 2).Value = 36.07
$ws2.Cells.Item(30, 3).Value = "historical"
$ws2.Cells.Item(31, 2).Value = 58.68331
$ws2.Cells.Item(31, 3).Value = "historical"
$ws2.Cells.Item(31, 5).Value = "OFF"
$ws2.Cells.Item(32, 2).Value = 57.06078
$ws2.Cells.Item(32, 3).Value = "historical"
$ws2.Cells.Item(33, 2).Value = 57.06035
$ws2.Cells.Item(33, 3).Value = "historical"
$ws2.Cells.Item(34, 2).Value = 68.88052
$ws2.Cells.Item(35, 2).Value = 62.65568
$ws2.Cells.Item(36, 2).Value = 66.82501000000001
$ws2.Cells.Item(37, 2).Value = 50.55353
$ws2.Cells.Item(38, 2).Value = 56.88407
$ws2.Cells.Item(39, 2).Value = 63.27083
$ws2.Cells.Item(40, 2).Value = 120.35934
$ws2.Cells.Item(41, 2).Value = 122.96668
$ws2.Cells.Item(42, 2).Value = 121.38263
$ws2.Cells.Item(43, 2).Value = 107.07418
$ws2.Cells.Item(44, 2).Value = 108.01
$ws2.Cells.Item(45, 2).Value = 108.89
$ws2.Cells.Item(49, 2).Value = 102.86014

# New rows 50-97 in Detailed sheet
$r = 50
$ws2.Cells.Item($r, 1).Value = 46072
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 84.79000000000001
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 51
$ws2.Cells.Item($r, 1).Value = 46072.02083333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 86.32680000000001
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 52
$ws2.Cells.Item($r, 1).Value = 46072.04166666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 85.65000000000001
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 53
$ws2.Cells.Item($r, 1).Value = 46072.0625
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 85.65000000000001
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 54
$ws2.Cells.Item($r, 1).Value = 46072.08333333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 85.65000000000001
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 55
$ws2.Cells.Item($r, 1).Value = 46072.10416666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 79.95041999999999
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 56
$ws2.Cells.Item($r, 1).Value = 46072.125
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 79.95041999999999
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 57
$ws2.Cells.Item($r, 1).Value = 46072.14583333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 79.95041999999999
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 58
$ws2.Cells.Item($r, 1).Value = 46072.16666666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 79.9504
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 59
$ws2.Cells.Item($r, 1).Value = 46072.1875
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 79.95
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 60
$ws2.Cells.Item($r, 1).Value = 46072.20833333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 84.79000000000001
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 61
$ws2.Cells.Item($r, 1).Value = 46072.22916666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 100.76432
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 62
$ws2.Cells.Item($r, 1).Value = 46072.25
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 108.89
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 63
$ws2.Cells.Item($r, 1).Value = 46072.27083333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 133.44919
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 64
$ws2.Cells.Item($r, 1).Value = 46072.29166666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 97.53870000000001
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 65
$ws2.Cells.Item($r, 1).Value = 46072.3125
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 57.06
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 66
$ws2.Cells.Item($r, 1).Value = 46072.33333333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 51.80245
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 67
$ws2.Cells.Item($r, 1).Value = 46072.35416666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 36.06
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 68
$ws2.Cells.Item($r, 1).Value = 46072.375
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 35.60254
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 69
$ws2.Cells.Item($r, 1).Value = 46072.39583333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 36.06
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 70
$ws2.Cells.Item($r, 1).Value = 46072.41666666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 30.123
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 71
$ws2.Cells.Item($r, 1).Value = 46072.4375
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 24.928
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 72
$ws2.Cells.Item($r, 1).Value = 46072.45833333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 0.51
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 73
$ws2.Cells.Item($r, 1).Value = 46072.47916666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 0.51
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 74
$ws2.Cells.Item($r, 1).Value = 46072.5
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 6.93964
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 75
$ws2.Cells.Item($r, 1).Value = 46072.52083333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 23.19082
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 76
$ws2.Cells.Item($r, 1).Value = 46072.54166666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 36.06
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 77
$ws2.Cells.Item($r, 1).Value = 46072.5625
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 35.86
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 78
$ws2.Cells.Item($r, 1).Value = 46072.58333333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 35.88
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 79
$ws2.Cells.Item($r, 1).Value = 46072.60416666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 54.31168
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "ON"
$r = 80
$ws2.Cells.Item($r, 1).Value = 46072.625
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 57.06033
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 81
$ws2.Cells.Item($r, 1).Value = 46072.64583333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 56.98
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 82
$ws2.Cells.Item($r, 1).Value = 46072.66666666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 64.89
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 83
$ws2.Cells.Item($r, 1).Value = 46072.6875
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 57.06
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 84
$ws2.Cells.Item($r, 1).Value = 46072.70833333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 62.89669
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 85
$ws2.Cells.Item($r, 1).Value = 46072.72916666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 55.92615
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 86
$ws2.Cells.Item($r, 1).Value = 46072.75
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 69.87746
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 87
$ws2.Cells.Item($r, 1).Value = 46072.77083333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 137.95
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 88
$ws2.Cells.Item($r, 1).Value = 46072.79166666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 237.97
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 89
$ws2.Cells.Item($r, 1).Value = 46072.8125
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 226.57464
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 90
$ws2.Cells.Item($r, 1).Value = 46072.83333333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 197.43404
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 91
$ws2.Cells.Item($r, 1).Value = 46072.85416666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 131.64387
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 92
$ws2.Cells.Item($r, 1).Value = 46072.875
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 102.14561
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 93
$ws2.Cells.Item($r, 1).Value = 46072.89583333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 103.9107
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 94
$ws2.Cells.Item($r, 1).Value = 46072.91666666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 78
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 95
$ws2.Cells.Item($r, 1).Value = 46072.9375
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 91.09668000000001
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 96
$ws2.Cells.Item($r, 1).Value = 46072.95833333334
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 78
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
$r = 97
$ws2.Cells.Item($r, 1).Value = 46072.97916666666
$ws2.Cells.Item($r, 1).NumberFormat = $ws2.Cells.Item(2, 1).NumberFormat
$ws2.Cells.Item($r, 2).Value = 84.79000000000001
$ws2.Cells.Item($r, 3).Value = "forecast"
$ws2.Cells.Item($r, 4).Value = 46072
$ws2.Cells.Item($r, 4).NumberFormat = $ws2.Cells.Item(2, 4).NumberFormat
$ws2.Cells.Item($r, 5).Value = "OFF"
